# Admin panel change: update product rows, add a new product row with a
# hyperlinked image URL, and widen the productImage column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column B (productImage) to fit the longer URL values.
$ws.Columns.Item(2).ColumnWidth = 26.5

# --- Update existing row 2 ("Rich Dad Poor Dad" -> "Gunhoa ka devta") ---
$ws.Range("A2").Value = "Gunhoa ka devta"
$ws.Range("C2").Value = 234
$ws.Range("D2").Value = "law Book"
$ws.Range("E2").Value = "Hindi"
$ws.Range("F2").Value = "this is book"
$ws.Range("G2").Value = $true

# B2 keeps the "Hyperlink" visual style but (per the target workbook) no
# longer carries a live hyperlink, so add then immediately remove the link
# while keeping the applied style/formatting.
$ws.Hyperlinks.Add($ws.Range("B2"), "https://images-na.ssl-images-amazon.com/images/S/compressed.photo.goodreads.com/books/1517755071i/3282557.jpg") | Out-Null
$ws.Range("B2").Value = "https://images-na.ssl-images-amazon.com/images/S/compressed.photo.goodreads.com/books/1517755071i/3282557.jpg"
$ws.Hyperlinks.Item(1).Delete()

# --- Add new row 3 for the second product ---
$ws.Range("A3").Value = "suraj ka satva ghoda"
$ws.Range("C3").Value = 455
$ws.Range("D3").Value = "Information Technology"
$ws.Range("E3").Value = "English"
$ws.Range("F3").Value = "this is book"
$ws.Range("G3").Value = $true

$ws.Hyperlinks.Add($ws.Range("B3"), "https://images-na.ssl-images-amazon.com/images/S/compressed.photo.goodreads.com/books/1357204311i/6720421.jpg") | Out-Null
$ws.Range("B3").Value = "https://images-na.ssl-images-amazon.com/images/S/compressed.photo.goodreads.com/books/1357204311i/6720421.jpg"

# Update the active selection to match the edited cell.
$ws.Range("E3").Select()
